$d = $word.ActiveDocument

function XmlEscape($text) {
    $t = $text -replace '&', '&amp;'
    $t = $t -replace '<', '&lt;'
    $t = $t -replace '>', '&gt;'
    return $t
}

$rPrBase = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rPrSuper = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr>'

function MakeRun($text, $rpr, $pageBreakBefore) {
    $escText = XmlEscape $text
    $pb = ''
    if ($pageBreakBefore) { $pb = '<w:lastRenderedPageBreak/>' }
    $needsPreserve = ($text.Length -gt 0) -and (($text.Substring(0,1) -match '\s') -or ($text.Substring($text.Length-1,1) -match '\s'))
    if ($needsPreserve) {
        $tOpen = '<w:t xml:space="preserve">'
    } else {
        $tOpen = '<w:t>'
    }
    return '<w:r>' + $rpr + $pb + $tOpen + $escText + '</w:t></w:r>'
}

function MakePackageXml($bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$pPrFirstLine = '<w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:contextualSpacing/>' + $rPrBase + '</w:pPr>'

# ---- Paragraph 1 (Jackson Palmer / Dogecoin creation story) ----
$p1Run1 = 'Dogecoin has one of the most amazing stories to tell in its creation. Way back in 2013, a few years after the explosion of cryptocurrencies, the internet was raving on a new internet meme called Doge. This created an amazing like renaissance of creativity on the internet and this sparked the creativity of one man in particular. Jackson Palmer in 2013 jokingly tweeted “Investing in Dogecoins, pretty sure it’s the next big thing”. While he was tweeting this, he was currently working for Adobe Systems in Sydney, Australia with a passion in cryptocurrency. Several people retweeted him and convinced him to purse his idea and not long after he purchased the domain name of dogecoin.com. although he did this for fun in the beginning it is no longer a joke. However, dogecoin still '
$p1Run2 = 'carries'
$p1Run3 = ' the title of “joke currency” to this day. But he was not alone in his '
$p1Run4 = 'endeavors'
$p1Run5 = '.'

$para1 = '<w:p>' + $pPrFirstLine `
    + (MakeRun $p1Run1 $rPrBase $true) `
    + (MakeRun $p1Run2 $rPrBase $false) `
    + (MakeRun $p1Run3 $rPrBase $false) `
    + (MakeRun $p1Run4 $rPrBase $false) `
    + (MakeRun $p1Run5 $rPrBase $false) `
    + '</w:p>'

# ---- Paragraph 2 (Billy Markus / launch details) ----
$p2Run1 = 'While Jackson was working on his crypto currency another person who was pursing his dreams in another cryptocurrency (which was not turning out so well) found the domain of dogecoin and reached out to Jackson. Billy Markus the second creator of dogecoin, lived in Portland reached out to Jackson and together they created dogecoin. Dogecoin officially launched on December 6'
$p2Run2 = 'th'
$p2Run3 = ' 2013, and the website initially saw over a million people on the site. With dogecoin it was initially limited at 100 billion dogecoin but soon changed to unlimited dogecoin.'
$p2Run4 = ' Seeing the great success of dogecoin the pair decide to do some marketing for the cryptocurrency and that is when reddit came in. The site was already trending and the site saw a market value of $8 million in the first two weeks. On of the key reasons that dogecoin was so successful is that unlike bitcoin people in the market of cryptocurrency already knew how crypto worked as there was less hesitancy. Another good reason that dogecoin took off was due to the bad background of bitcoin and the new laws set in place like in China banning the use of bitcoin, '
$p2Run5 = 'boosting the sales of dogecoin even more.'
$p2Run6 = ' Dogecoin today is a top contender in the cryptocurrency game and its still looking to be a popular choice for most and even Elon Musk. '

$para2 = '<w:p>' + $pPrFirstLine `
    + (MakeRun $p2Run1 $rPrBase $false) `
    + (MakeRun $p2Run2 $rPrSuper $false) `
    + (MakeRun $p2Run3 $rPrBase $false) `
    + (MakeRun $p2Run4 $rPrBase $false) `
    + (MakeRun $p2Run5 $rPrBase $true) `
    + (MakeRun $p2Run6 $rPrBase $false) `
    + '</w:p>'

# ---- Insert both paragraphs after the current last paragraph ----
$lastPara = $d.Paragraphs.Last
$endRng = $lastPara.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

$target = $d.Paragraphs.Item(9).Range
$target.Collapse(1)
$target.InsertXML((MakePackageXml ($para1 + $para2)))

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
